$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("17:55", "Epiphanie Adoboè", "CEO", "Départ"),
  @("17:40", "Akouété Kangnivi", "Juriste", "Départ"),
  @("17:35", "Ibrahim Kodjo", "Web", "Départ"),
  @("16:51", "Amah Kwatcha", "CEO", "Départ"),
  @("16:45", "Abi Conrad", "CEO", "Départ"),
  @("16:45", "Abalo Afi", "CEO", "Départ"),
  @("16:38", "Romuald Gagnon", "PHOTOGRAPHE", "Départ"),
  @("16:15", "Adjovi Abla", "Developpeur", "Départ"),
  @("13:51", "Amah Kwatcha", "CEO", "Fin pause"),
  @("13:45", "Abi Conrad", "CEO", "Fin pause"),
  @("13:22", "Napo Kuvor", "Mecanicien", "Fin pause"),
  @("11:51", "Amah Kwatcha", "CEO", "Pause"),
  @("11:45", "Abi Conrad", "CEO", "Pause"),
  @("11:22", "Napo Kuvor", "Mecanicien", "Pause"),
  @("08:55", "Epiphanie Adoboè", "CEO", "Arrivée"),
  @("08:40", "Akouété Kangnivi", "Juriste", "Arrivée"),
  @("08:35", "Ibrahim Kodjo", "Web", "Arrivée"),
  @("07:51", "Amah Kwatcha", "CEO", "Arrivée"),
  @("07:45", "Abi Conrad", "CEO", "Arrivée"),
  @("07:45", "Abalo Afi", "CEO", "Arrivée"),
  @("07:38", "Romuald Gagnon", "PHOTOGRAPHE", "Arrivée"),
  @("07:22", "Napo Kuvor", "Mecanicien", "Arrivée"),
  @("07:15", "Adjovi Abla", "Developpeur", "Arrivée")
)

$r = 4
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

$ws.Rows.Item(27).Delete()
